$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.755.74'
$ws.Range('E2').Value = '  -2.40%  '
$ws.Range('D3').Value = '1.565.93'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.31'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.90%  '
$ws.Range('E6').Value = '  -2.26%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '21.93'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.26%  '
$ws.Range('E9').Value = '  -0.41%  '
$ws.Range('E10').Value = '  -1.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0862'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.49%  '
$ws.Range('D12').Value = '1.788.19'
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D13').Value = '1.557.52'
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('E14').Value = '  -2.38%  '
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('D16').Value = '26.803.51'
$ws.Range('E16').Value = '  -2.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.51'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.07%  '
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '214.37'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.96%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.35'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.57%  '
$ws.Range('D20').Value = '0.0₃0677'
$ws.Range('E20').Value = '  -1.62%  '
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.09'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('E23').Value = '  -1.69%  '
$ws.Range('E24').Value = '  -1.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.78'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('E26').Value = '  +0.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.91'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('E29').Value = '  -1.36%  '
$ws.Range('E30').Value = '  -1.32%  '
$ws.Range('E31').Value = '  -3.65%  '
$ws.Range('E32').Value = '  -1.45%  '
$ws.Range('D33').Value = '1.386.13'
$ws.Range('E33').Value = '  +1.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.93'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.98%  '
$ws.Range('E35').Value = '  +0.94%  '
$ws.Range('E36').Value = '  -0.74%  '
$ws.Range('E37').Value = '  -3.99%  '
$ws.Range('E38').Value = '  -2.46%  '
$ws.Range('E39').Value = '  -1.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.816'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.35%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.989'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.45%  '
$ws.Range('E43').Value = '  +1.85%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.78'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.83%  '
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.18'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.66%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '63.29'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.71%  '
$ws.Range('D47').Value = '1.701.31'
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '85.42'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.16%  '
$ws.Range('D49').Value = '0.0₇0985'
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('E50').Value = '  -0.44%  '
$ws.Range('E51').Value = '  -0.70%  '
